# Adds "(рис. N)" references to the 16 procedure-step paragraphs and, after
# each one, inserts a blank paragraph followed by an italic figure caption
# paragraph ("Рис. N: ...").

function New-FlatOpc($bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyInner + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Replaces the whole paragraph's content with a single plain-text run
# (this also clears any pre-existing trailing runs, e.g. the lone " " run).
function Set-ParagraphPlainText($para, $text) {
    $inner = '<w:body><w:p><w:r><w:t xml:space="preserve">' + $text + '</w:t></w:r></w:p></w:body>'
    $para.Range.InsertXML((New-FlatOpc $inner))
}

# Inserts a bare empty paragraph (no runs) right after $para; returns it.
function Add-BlankParagraphAfter($para) {
    $para.Range.InsertParagraphAfter()
    $newPara = $para.Next()
    $inner = '<w:body><w:p/></w:body>'
    $newPara.Range.InsertXML((New-FlatOpc $inner))
    return $newPara
}

# Inserts an italic caption paragraph right after $para; returns it.
function Add-ItalicCaptionAfter($para, $text) {
    $para.Range.InsertParagraphAfter()
    $newPara = $para.Next()
    $inner = '<w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">' + $text + '</w:t></w:r></w:p></w:body>'
    $newPara.Range.InsertXML((New-FlatOpc $inner))
    return $newPara
}

$d = $word.ActiveDocument

$items = @(
  @("1: Определяю имя домашнего каталога (рис. 1)", "Рис. 1: Определение имени домашнего каталога"),
  @("2.1: Перешел в tmp, показал содержимое без аргументов... (рис. 2.1)", "Рис. 2.1: Содержимое tmp без аргументов"),
  @("2.2.2: ...С аргументом a... (рис. 2.2)", "Рис. 2.2: Содержимое tmp с аргументом a"),
  @("2.2.3: ...С аргументом F (рис. 2.3)", "Рис. 2.3: Содержимое tmp с аргументом F"),
  @("2.2.4: ... И с аргументом l (рис. 2.4)", "Рис. 2.4: Содержимое tmp с аргументом l"),
  @("2.3-2.4: cron есть в var/spoon, владелец файлов - я, ksanikin (рис. 2.5)", "Рис. 2.5: cron и ksanikin"),
  @("3.1-3.3: Создаю newdir и /morefun, создаю и удаляю три папки сразу (рис. 3.1)", "Рис. 3.1: Работа с mkdir и rmdir"),
  @("3.4-3.5: Удаляю каталоги newdir и morefun, rm с newdir не работает (рис. 3.2)", "Рис. 3.2: Работа с mkdir и rmdir"),
  @("4: С помощью mana нашел аргумент для выведения подкаталогов - R (рис. 4)", "Рис. 4: Выведение подкаталогов"),
  @("5: С помощью mana нашел аргумент для сортировки по времени - t Аргумент для выведения описания файлов - l Полный аргумент - lt (рис. 5)", "Рис. 5: man ls"),
  @("6cd: -L следует по ссылкам на файлы -P не следует -eP выдает ошибку, если cd не может сменить папку (рис. 6.1)", "Рис. 6.1: man cd"),
  @("6pwd: -L указывает каталог с учетом ссылок -P без учета (рис. 6.2)", "Рис. 6.2: man pwd"),
  @("6mkdir: -m указывает права доступа папки -p создаст родительские папки в структуре папок при их отсутствии -v выведет сообщение о каждой созданной папке -Z ставит контекст безопасности папок по умолчанию, что бы это ни значило (рис. 6.3)", "Рис. 6.3: man mkdir"),
  @("6rmdir: -p удаляет каталог вместе с файлами -v выведет сообщение о каждой удаленной папке -i при удалении каждой будет выведен запрос на подтверждение (рис. 6.4)", "Рис. 6.4: man rmdir"),
  @("6rm: -f не выводит ошибку при удалении несуществующих файлов -i при удалении каждого будет выведен запрос на подтверждение -d аналог rmdir -v выведет сообщение о каждой удаленном удаленном файле (рис. 6.5)", "Рис. 6.5: man rm"),
  @("7: Воспользовался history, запустил программу из истории с измененным параметром (рис. 7)", "Рис. 7: Работа с history")
)

# The 16 target paragraphs start at index 9 in the original document and,
# since each processed item leaves the *next* original paragraph exactly
# 3 slots further along (itself + blank + caption), a simple running index
# keeps everything aligned without needing to search for text each time.
$idx = 9
foreach ($pair in $items) {
    $newText = $pair[0]
    $capText = $pair[1]

    $p = $d.Paragraphs($idx)
    Set-ParagraphPlainText $p $newText

    $p = $d.Paragraphs($idx)
    $blank = Add-BlankParagraphAfter $p
    Add-ItalicCaptionAfter $blank $capText | Out-Null

    $idx = $idx + 3
}

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
